$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell (existing last data row) used to copy the date number format
$dateFormat = $ws.Cells.Item(248, 4).NumberFormat

$rows = @(
    @{ A=8; B="Terminal La Palmera de La Serena"; C="Coquimbo"; D=44595; E=4; F="Fruta"; G=100103; H="Frutos de hueso (carozo)"; I=100103006; J="Nectarín"; K="June Pearl";  L="Especial"; M=20; N=400000; O=405000; P=402500; Q="$/bins (420 kilos)"; R="Región de O'Higgins"; S=958; T=420 },
    @{ A=8; B="Terminal La Palmera de La Serena"; C="Coquimbo"; D=44595; E=4; F="Fruta"; G=100103; H="Frutos de hueso (carozo)"; I=100103006; J="Nectarín"; K="June Pearl";  L="Primera";  M=20; N=350000; O=360000; P=355000; Q="$/bins (420 kilos)"; R="Región de O'Higgins"; S=845; T=420 },
    @{ A=8; B="Terminal La Palmera de La Serena"; C="Coquimbo"; D=44595; E=4; F="Fruta"; G=100103; H="Frutos de hueso (carozo)"; I=100103006; J="Nectarín"; K="June Pearl";  L="Segunda";  M=20; N=325000; O=330000; P=327500; Q="$/bins (420 kilos)"; R="Región de O'Higgins"; S=780; T=420 },
    @{ A=8; B="Terminal La Palmera de La Serena"; C="Coquimbo"; D=44595; E=4; F="Fruta"; G=100103; H="Frutos de hueso (carozo)"; I=100103006; J="Nectarín"; K="Venus";      L="Especial"; M=20; N=370000; O=375000; P=372500; Q="$/bins (420 kilos)"; R="Región de O'Higgins"; S=887; T=420 },
    @{ A=8; B="Terminal La Palmera de La Serena"; C="Coquimbo"; D=44595; E=4; F="Fruta"; G=100103; H="Frutos de hueso (carozo)"; I=100103006; J="Nectarín"; K="Venus";      L="Primera";  M=28; N=330000; O=335000; P=331786; Q="$/bins (420 kilos)"; R="Región de O'Higgins"; S=790; T=420 }
)

$startRow = 249
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data.A
    $ws.Cells.Item($r, 2).Value = $data.B
    $ws.Cells.Item($r, 3).Value = $data.C

    $ws.Cells.Item($r, 4).Value = $data.D
    $ws.Cells.Item($r, 4).NumberFormat = $dateFormat

    $ws.Cells.Item($r, 5).Value = $data.E
    $ws.Cells.Item($r, 6).Value = $data.F
    $ws.Cells.Item($r, 7).Value = $data.G
    $ws.Cells.Item($r, 8).Value = $data.H
    $ws.Cells.Item($r, 9).Value = $data.I
    $ws.Cells.Item($r, 10).Value = $data.J
    $ws.Cells.Item($r, 11).Value = $data.K
    $ws.Cells.Item($r, 12).Value = $data.L
    $ws.Cells.Item($r, 13).Value = $data.M
    $ws.Cells.Item($r, 14).Value = $data.N
    $ws.Cells.Item($r, 15).Value = $data.O
    $ws.Cells.Item($r, 16).Value = $data.P
    $ws.Cells.Item($r, 17).Value = $data.Q
    $ws.Cells.Item($r, 18).Value = $data.R
    $ws.Cells.Item($r, 19).Value = $data.S
    $ws.Cells.Item($r, 20).Value = $data.T
}

Write-Host "Added rows 249-253"
